$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("M2").Value = 370
$ws.Range("L3").Value = 7111
$ws.Range("M3").Value = 410
$ws.Range("B4").Value = 1724
$ws.Range("M4").Value = 113
$ws.Range("L5").Value = 424
$ws.Range("M5").Value = 25
$ws.Range("L6").Value = 5819
$ws.Range("M6").Value = 312
$ws.Range("B7").Value = 23356
$ws.Range("L7").Value = 21732
$ws.Range("M7").Value = 1230

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("M7").Value = 39
$ws.Range("M8").Value = 84
$ws.Range("M11").Value = 14
$ws.Range("M15").Value = 12
$ws.Range("M19").Value = 45
$ws.Range("M21").Value = 5
$ws.Range("M27").Value = 17
$ws.Range("M29").Value = 61
$ws.Range("M31").Value = 18
$ws.Range("M33").Value = 45
$ws.Range("L36").Value = 271
$ws.Range("M37").Value = 56
$ws.Range("M42").Value = 41
$ws.Range("M46").Value = 4
$ws.Range("M51").Value = 19
$ws.Range("M52").Value = 16
$ws.Range("M57").Value = 2
$ws.Range("M60").Value = 11
$ws.Range("B63").Value = 428
$ws.Range("L63").Value = 66
$ws.Range("M63").Value = 2
$ws.Range("M65").Value = 27
$ws.Range("M67").Value = 35
$ws.Range("M72").Value = 8
$ws.Range("L78").Value = 285
$ws.Range("M79").Value = 27
$ws.Range("M85").Value = 61
$ws.Range("M86").Value = 9
$ws.Range("M89").Value = 17
$ws.Range("M94").Value = 16
$ws.Range("M95").Value = 16
$ws.Range("M99").Value = 33
$ws.Range("B101").Value = 23356
$ws.Range("L101").Value = 21732
$ws.Range("M101").Value = 1230

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("M5").Value = 1
$ws.Range("M7").Value = 39

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("M3").Value = 3
$ws.Range("M7").Value = 14

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("M3").Value = 5
$ws.Range("M7").Value = 17

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("M3").Value = 27
$ws.Range("M6").Value = 12
$ws.Range("M7").Value = 61

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("M6").Value = 5
$ws.Range("M7").Value = 16

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("M2").Value = 24
$ws.Range("M3").Value = 28
$ws.Range("M6").Value = 26
$ws.Range("M7").Value = 84

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("M3").Value = 15
$ws.Range("M7").Value = 45

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("M6").Value = 4
$ws.Range("M7").Value = 16

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("M6").Value = 12
$ws.Range("M7").Value = 56

$ws = $wb.Worksheets.Item('New City')
$ws.Range("M2").Value = 8
$ws.Range("M3").Value = 11
$ws.Range("M4").Value = 2
$ws.Range("M7").Value = 27

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("M2").Value = 12
$ws.Range("M7").Value = 33

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("M6").Value = 5
$ws.Range("M7").Value = 18

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 214
$ws.Range("L3").Value = 295
$ws.Range("M3").Value = 11
$ws.Range("M7").Value = 35

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("M4").Value = 6
$ws.Range("M7").Value = 61

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("M6").Value = 13
$ws.Range("M7").Value = 45

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("M6").Value = 11
$ws.Range("M7").Value = 41

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L2").Value = 75
$ws.Range("L7").Value = 285

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("M6").Value = 1
$ws.Range("M7").Value = 4

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("M2").Value = 2
$ws.Range("M7").Value = 5

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("M6").Value = 6
$ws.Range("M7").Value = 27

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L6").Value = 66
$ws.Range("L7").Value = 271

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("M6").Value = 6
$ws.Range("M7").Value = 16

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("M3").Value = 4
$ws.Range("M7").Value = 12

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("M2").Value = 5
$ws.Range("M3").Value = 6
$ws.Range("M7").Value = 17

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("M4").Value = 6
$ws.Range("M7").Value = 9

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("M2").Value = 8
$ws.Range("M3").Value = 4
$ws.Range("M6").Value = 3
$ws.Range("M7").Value = 19

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("M3").Value = 2
$ws.Range("M7").Value = 2

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("M2").Value = 4
$ws.Range("M6").Value = 3
$ws.Range("M7").Value = 11

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("M6").Value = 3
$ws.Range("M7").Value = 8
